$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was inserted as row 380 ("Cuatro cascos verde" on
# 2023-01-20), pushing the existing rows 380-432 down to 381-433.
$ws.Rows.Item(380).Insert()

$ws.Cells.Item(380, 1).Value = 11
$ws.Cells.Item(380, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(380, 3).Value = "Bíobío"
$ws.Cells.Item(380, 4).Value = 44946
$ws.Cells.Item(380, 5).Value = 8
$ws.Cells.Item(380, 6).Value = 100112002
$ws.Cells.Item(380, 7).Value = "Pimiento"
$ws.Cells.Item(380, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(380, 9).Value = "Primera"
$ws.Cells.Item(380, 10).Value = 180
$ws.Cells.Item(380, 11).Value = 11000
$ws.Cells.Item(380, 12).Value = 12000
$ws.Cells.Item(380, 13).Value = 11556
$ws.Cells.Item(380, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(380, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(380, 16).Value = 642
$ws.Cells.Item(380, 17).Value = 18
$ws.Cells.Item(380, 18).Value = "Hortaliza"
